$wb = $excel.ActiveWorkbook

# The edit applies to the second worksheet, "Tasks 01-28 to 02-04"
$ws = $wb.Worksheets.Item("Tasks 01-28 to 02-04")

# Update the note on row 8 (Create ComparedItem Game Object) - the old note about
# comparing equipped item info is replaced by a new note about magic comparison.
$ws.Range("G8").Value = "Still need to alter to show compared magic when looking at magic section"

# Add a new task row (row 14): "Create Icon Base Class"
$ws.Range("A14").Value = "Create Icon Base Class"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = 2
$ws.Range("E14").Value = "James"
$ws.Range("F14").Value = "Done"
$ws.Range("F14").Interior.Color = 5296274
$ws.Range("G14").Value = "Icon Spawner can spawn these in"

# Recalculate totals (B16/C16 SUM formulas already cover rows up to 14/15)
$excel.Calculate()

# Update the selected/active cell to reflect the last edited cell
$ws.Activate() | Out-Null
$ws.Range("G14").Select() | Out-Null
